# Re-order the header columns of the "Products" and "Sales" sheets and add
# a new data row (row 2) to each, reflecting a single product + its sale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Products" (sheet1)
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")

$productsHeaders = @{
    "A1" = "name";
    "B1" = "description";
    "C1" = "price";
    "D1" = "cost_price";
    "E1" = "category";
    "F1" = "stock";
    "G1" = "min_stock";
    "H1" = "supplier";
    "I1" = "sku";
    "J1" = "id";
    "K1" = "created_date";
    "L1" = "last_updated";
}
foreach ($addr in $productsHeaders.Keys) {
    $wsProducts.Range($addr).Value = $productsHeaders[$addr]
}

# Text-valued cells for the new data row (forced to text so numeric-looking
# strings such as "100" / a GUID / an ISO date are not coerced to numbers).
$productsRow2Text = @{
    "A2" = "amul";
    "B2" = "";
    "C2" = "100";
    "D2" = "50";
    "E2" = "Food";
    "H2" = "";
    "I2" = "";
    "J2" = "285f3854-cbf2-4ba5-b477-5cda3035ab2d";
    "K2" = "2025-09-12T12:06:21.787Z";
    "L2" = "2025-09-12T12:06:48.292Z";
}
foreach ($addr in $productsRow2Text.Keys) {
    $wsProducts.Range($addr).NumberFormat = "@"
    $wsProducts.Range($addr).Value = $productsRow2Text[$addr]
}

# Numeric cells for the new data row.
$wsProducts.Range("F2").Value = 80
$wsProducts.Range("G2").Value = 1

# ---------------------------------------------------------------------
# Sheet "Sales" (sheet2)
# ---------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Sales")

$salesHeaders = @{
    "A1" = "product_id";
    "B1" = "quantity";
    "C1" = "unit_price";
    "D1" = "customer_name";
    "E1" = "payment_method";
    "F1" = "cashier";
    "G1" = "notes";
    "H1" = "id";
    "I1" = "product_name";
    "J1" = "total_amount";
    "K1" = "profit";
    "L1" = "sale_date";
}
foreach ($addr in $salesHeaders.Keys) {
    $wsSales.Range($addr).Value = $salesHeaders[$addr]
}

$salesRow2Text = @{
    "A2" = "285f3854-cbf2-4ba5-b477-5cda3035ab2d";
    "C2" = "100";
    "D2" = "Zeeshan";
    "E2" = "Cash";
    "F2" = "Admin";
    "G2" = "";
    "H2" = "937b64d0-f047-42cc-b7ea-14a294d8ed1c";
    "I2" = "amul";
    "J2" = "2000";
    "K2" = "1000";
    "L2" = "2025-09-12T12:06:48.282Z";
}
foreach ($addr in $salesRow2Text.Keys) {
    $wsSales.Range($addr).NumberFormat = "@"
    $wsSales.Range($addr).Value = $salesRow2Text[$addr]
}

# Numeric cell for the new data row.
$wsSales.Range("B2").Value = 20
